# Fixed BoM and PnP files for CH340G and FT232RL boards.
#
# The workbook originally ships with an empty "Sheet1" plus a "Sheet2"
# that carries the actual Bill-of-Materials query table (and the
# ExternalData_1 defined name that points at it). The fix removes the
# stray empty sheet and promotes the data sheet to "Sheet1".

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Drop the empty placeholder sheet.
$wb.Worksheets.Item("Sheet1").Delete() | Out-Null

# Promote the remaining (data) sheet so it takes over the "Sheet1" name.
$wb.Worksheets.Item("Sheet2").Name = "Sheet1"

# Make sure the surviving sheet is the active one (it also was the
# active sheet in the original workbook).
$wb.Worksheets.Item("Sheet1").Select()
